$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Admin_Quiz" worksheet right before the existing "User"
#    sheet (Sheet1, Admin_Quiz, User).
# ---------------------------------------------------------------------------
$userWs = $wb.Worksheets.Item("User")
$quiz = $wb.Worksheets.Add($userWs)
$quiz.Name = "Admin_Quiz"

# ---------------------------------------------------------------------------
# 2. Populate Admin_Quiz with the new quiz-related API test data.
# ---------------------------------------------------------------------------

# Header row
$quiz.Range("A1").Value = "testcase_name"
$quiz.Range("B1").Value = "userName"
$quiz.Range("C1").Value = "firstName"
$quiz.Range("D1").Value = "lastName"
$quiz.Range("E1").Value = "password"
$quiz.Range("F1").Value = "email"
$quiz.Range("G1").Value = "phone"
$quiz.Range("H1").Value = "huBatchId"
$quiz.Range("I1").Value = "title"
$quiz.Range("J1").Value = "categoryId"
$quiz.Range("K1").Value = "question_number"
$quiz.Range("L1").Value = "huBatchTitle"
$quiz.Range("M1").Value = "startDate"
$quiz.Range("N1").Value = "endDate"

# Row 2
$quiz.Range("A2").Value = "add_single_quiz"
$quiz.Range("J2").Value = "13"

# Row 3
$quiz.Range("A3").Value = "add_category"
$quiz.Range("I3").Value = "Java Core"

# Row 4
$quiz.Range("A4").Value = "get_all_questions"
$quiz.Range("K4").Value = 9

# Row 5
$quiz.Range("A5").Value = "update_a_question"

# Row 6
$quiz.Range("A6").Value = "delete_a_question"
$quiz.Range("K6").Value = 1

# Row 7 (tall row for a long question body)
$quiz.Range("A7").Value = "add_question"
$quiz.Rows.Item(7).RowHeight = 82
$quiz.Range("L7").WrapText = $true

# Row 8
$quiz.Range("A8").Value = "save_hux"
$quiz.Range("L8").Value = "'=RANDBETWEEN(10,100)"
$quiz.Range("M8").Value = "2022-01-10"
$quiz.Range("N8").Value = "2022-02-14"

# Row 9
$quiz.Range("A9").Value = "get_all_hux"
$quiz.Range("L9").Value = "'=RANDBETWEEN(10,100)"
$quiz.Range("M9").Value = "2022-01-10"
$quiz.Range("N9").Value = "2022-02-14"

# Row 10
$quiz.Range("A10").Value = "get_alluser_by_huid"
$quiz.Range("H10").Value = 3

# Row 11 (login credentials + hyperlink)
$quiz.Range("A11").Value = "login"
$quiz.Range("B11").Value = "admin"
$quiz.Range("E11").Value = "1234@admin"
$quiz.Hyperlinks.Add($quiz.Range("E11"), "mailto:admin@test.com") | Out-Null

# Column widths
$quiz.Columns.Item(1).ColumnWidth = 20.5
$quiz.Columns.Item(11).ColumnWidth = 17.5
$quiz.Columns.Item(12).ColumnWidth = 65.6640625

# Selection / view state used by the authored workbook
$quiz.Range("F13").Select()

# ---------------------------------------------------------------------------
# 3. Extend the "User" sheet with the quiz title/category columns.
# ---------------------------------------------------------------------------
$userWs.Range("I1").Value = "title"
$userWs.Range("J1").Value = "categoryId"
$userWs.Range("A3").Value = "add_category"
$userWs.Range("J3").Value = "1"

$userWs.Range("A1:XFD3").Select()
